# The original "Test Name" header cell (A1) is renamed to "TestName", and
# a new first data row is inserted with the value "TestMethod1" in A2,
# shifting the old data rows down. The active selection moves to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TestName"
$ws.Range("A2").Value = "TestMethod1"

$ws.Range("A2").Select() | Out-Null
